# Correct text book references on slide 2
#
# Slide 2 ("Resources") has a body placeholder shape listing text-book
# section references. Two of the references were wrong and are corrected:
#   "Sec 5.1,5.2,5.4,5.8,5.9 "  ->  "Sec 7.1,7.2,7.3,7.4,7.5,8.2,11.1"
#   "Sec 8.2-8.4 "              ->  "Sec 12.1,12.2"
#
# The "7." prefix of the first corrected reference is emphasised (bold +
# single underline) while keeping the Courier New font already used for
# the whole reference run. The previously-trailing " " run in both
# paragraphs is removed (the replacement text no longer needs it).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# --- Paragraph "Sec 5.1,5.2,5.4,5.8,5.9 " (Text book 2: Horowitz) ---
# The Courier-New run holding "5.1,5.2,5.4,5.8,5.9" starts right after
# "Sec " at character 27 and is 19 characters long.
$oldRef1 = $tr.Characters(27, 19)

# Prepend "7." (inherits the surrounding Courier New run's formatting).
$oldRef1.InsertBefore("7.") | Out-Null

# Make the newly inserted "7." bold + underlined - this splits it into
# its own run, distinct from the (non-bold) remainder of the reference.
$newPrefix = $tr.Characters(27, 2)
$newPrefix.Font.Bold = 1
$newPrefix.Font.Underline = 1

# Replace the remaining original digits "5.1,5.2,5.4,5.8,5.9" (now
# shifted two characters later, at 29) with the corrected reference list.
$remainder1 = $tr.Characters(29, 19)
$remainder1.Text = "1,7.2,7.3,7.4,7.5,8.2,11.1"

# Drop the trailing " " run that used to follow the reference.
$trailingSpace1 = $tr.Characters(55, 1)
$trailingSpace1.Text = ""

# --- Paragraph "Sec 8.2-8.4 " (Text book 1: Levitin) ---
# The Courier-New run holding "8.2-8.4" starts at character 81 and is 7
# characters long.
$oldRef2 = $tr.Characters(81, 7)
$oldRef2.Text = "12.1,12.2"

# Drop the trailing " " run that used to follow the reference (now at 90
# since the replacement text is two characters longer than the original).
$trailingSpace2 = $tr.Characters(90, 1)
$trailingSpace2.Text = ""
